$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new record row at row 297 (pushes existing rows 297:331 down to 298:332)
$ws.Rows.Item(297).Insert()

$ws.Cells.Item(297, 1).Value = 10
$ws.Cells.Item(297, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(297, 3).Value = "La Araucanía"
$ws.Cells.Item(297, 4).Value = 45034
$ws.Cells.Item(297, 5).Value = 9
$ws.Cells.Item(297, 6).Value = 100112043
$ws.Cells.Item(297, 7).Value = "Pepino dulce"
$ws.Cells.Item(297, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(297, 9).Value = "Primera"
$ws.Cells.Item(297, 10).Value = 55
$ws.Cells.Item(297, 11).Value = 16000
$ws.Cells.Item(297, 12).Value = 17000
$ws.Cells.Item(297, 13).Value = 16545
$ws.Cells.Item(297, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(297, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(297, 16).Value = 919
$ws.Cells.Item(297, 17).Value = 18
$ws.Cells.Item(297, 18).Value = "Hortaliza"
